$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49-74 down to 50-75
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 45240
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100112026
$ws.Range("G49").Value = "Haba"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 50
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = 9000
$ws.Range("N49").Value = "$/saco 25 kilos"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 360
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"

# Make sure the date cell keeps the same number format as the rest of column D
$ws.Range("D49").NumberFormat = $ws.Range("D50").NumberFormat()
